$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the D2 and D5 values (165 <-> 161); formulas in E2/E5 recalc automatically
$ws.Range("D2").Value = 161
$ws.Range("D5").Value = 165

# Update the selected cell/range from G5 to F5
$ws.Range("F5").Select()
